$d = $word.ActiveDocument

# The heading currently reads "OFÍCIO/IDAF/SIMLAM Nº «Numero»", built out of
# the runs: [OFÍCIO] [_GoBack bookmark] [/IDAF] [/SIMLAM] [ N] [º ] [Numero field].
#
# Target: "OFÍCIO/IDAF Nº «Numero»", i.e. drop the "/SIMLAM" run (and its
# neighboring old "/IDAF" run) and re-insert "/IDAF" right after "OFÍCIO",
# ahead of the _GoBack bookmark, leaving the bookmark marking the same spot
# between the (new) "/IDAF" and the following " N" run.

# 1) Remove the old "/IDAF/SIMLAM" text that sits right after the bookmark.
$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute("/IDAF/SIMLAM", $true, $false, $false, $false, $false, $true, 1, $false)
if ($found) {
    $rng.Delete()
}

# 2) Re-insert "/IDAF" immediately after "OFÍCIO", before the bookmark.
$find = $d.Content
$find.Find.ClearFormatting()
$found2 = $find.Find.Execute("OFÍCIO", $true, $false, $false, $false, $false, $true, 1, $false)
if ($found2) {
    $insertPoint = $d.Range($find.End, $find.End)
    $insertPoint.InsertAfter("/IDAF")
}

$d.Save()
